# Daily attendance processing - 2026-01-16 23:03:14
# For every row in the "Recorded By" column (G), if the comma-separated
# list of recorders contains an exact "System" entry, reverse the order
# of the whole list (e.g. "System, a, b" -> "b, a, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $raw = $cell.Value2

    if ($raw -eq $null) { continue }

    $value = [string]$raw
    if ($value -eq "") { continue }

    $parts = $value.Split(",")
    $trimmedParts = @()
    foreach ($part in $parts) {
        $trimmedParts += $part.Trim()
    }

    $hasSystem = $false
    foreach ($part in $trimmedParts) {
        if ($part -ceq "System") {
            $hasSystem = $true
        }
    }

    if ($hasSystem -and $trimmedParts.Count -gt 1) {
        $reversed = $trimmedParts[($trimmedParts.Count - 1)..0]
        $newValue = $reversed -join ", "
        $cell.Value = $newValue
    }
}
